$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 2-36) held the counts as text (shared-string) cells, e.g.
# B2 = "237" (text). Re-enter them as real numbers so the cells store a
# numeric <v> instead of a shared-string reference. The last row (B36)
# also changes value from 1 to 0 (the previous "1" next to the bogus
# "err" label was wrong).
$counts = @(237,230,213,204,202,190,186,169,159,158,149,147,147,144,143,136,131,130,129,126,124,116,105,92,65,49,40,31,25,20,15,4,2,1,0)
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

# A36 previously held the literal text "err" (a leftover from a failed
# lookup). Fix it to read "-1.0" while keeping it a plain text cell (not
# a number) and without touching the cell's style. Enter it as a formula
# that yields the text, then flatten it to a static value via
# copy/paste-special so no formula or style residue is left behind.
$ws.Range("A36").Formula = '="-1.0"'
$ws.Range("A36").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# Make sure the chart picks up the corrected category label / counts.
if ($ws.ChartObjects().Count -gt 0) {
    $ws.ChartObjects().Item(1).Chart.Refresh()
}
